# Add a new "RESTAURANT USERS" table block (a login-related user table)
# to the database-tables overview sheet, in column I starting at row 25,
# mirroring the other table blocks already present on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I25").Value = "RESTAURANT USERS"
$ws.Range("I26").Value = "user_id"
$ws.Range("I27").Value = "email"
$ws.Range("I28").Value = "password"
$ws.Range("I29").Value = "account_type"
$ws.Range("I30").Value = "name"
$ws.Range("I31").Value = "surname"
$ws.Range("I32").Value = "restaurant_id"

# Match the author's final cursor position/selection on the sheet.
$ws.Activate() | Out-Null
$ws.Range("I30").Select() | Out-Null
